$d = $word.ActiveDocument

# 1) "English" -> Arabic (both occurrences: hyperlink + plain text)
$d.Content.Find.Execute("English", $true, $false, $false, $false, $false, $true, 1, $false, "الإنجليزية", 2) | Out-Null

# 2) " / Portuguese / French / Thai / Vietnamese / Spanish" -> Arabic language list
$d.Content.Find.Execute(" / Portuguese / French / Thai / Vietnamese / Spanish", $true, $false, $false, $false, $false, $true, 1, $false, " /البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية", 2) | Out-Null

# 3) "Brief" -> "المضمون"
$d.Content.Find.Execute("Brief", $true, $false, $false, $false, $false, $true, 1, $false, "المضمون", 2) | Out-Null

# 4) Brief description sentence (partial -> keep customer.io, translate tail)
$d.Content.Find.Execute("An email sent to partners in the target country who RSVPed yes but haven’t sent their documents to us. It will be sent via customer.io", $true, $false, $false, $false, $false, $true, 1, $false, "An email sent to partners in the target country who RSVPed yes but haven’t sent their documents to us. سيتم إرسالها عبر customer.io", 2) | Out-Null

# 5) "Target audience" -> "الجمهور المستهدف"
$d.Content.Find.Execute("Target audience", $true, $false, $false, $false, $false, $true, 1, $false, "الجمهور المستهدف", 2) | Out-Null

# 6) "Don't forget to send your documents" -> Arabic (both occurrences)
$d.Content.Find.Execute("Don’t forget to send your documents", $true, $false, $false, $false, $false, $true, 1, $false, "لا تنس إرسال مستنداتك", 2) | Out-Null

# 7) "Hi " -> "مرحبًا  " (greeting run)
$d.Content.Find.Execute("Hi ", $true, $false, $false, $false, $false, $true, 1, $false, "مرحبًا  ", 2) | Out-Null

# 8) ", " -> ",، " -- ONLY in the "Hi [PARTNER NAME], " paragraph, not the "Dear" one.
#    Scope the Find to that specific paragraph (the one that now starts with the Arabic greeting).
$greetingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptxt = $d.Paragraphs.Item($i).Range.Text
    if ($ptxt.StartsWith("مرحبًا  ")) {
        $greetingPara = $d.Paragraphs.Item($i)
        break
    }
}
$greetingPara.Range.Find.Execute(", ", $true, $false, $false, $false, $false, $true, 1, $false, ",، ", 2) | Out-Null

# 9) "If you have any questions, please contact your country manager." -> Arabic
$d.Content.Find.Execute("If you have any questions, please contact your country manager.", $true, $false, $false, $false, $false, $true, 1, $false, "إذا كانت لديك أي أسئلة، فيُرجى الاتصال بمديرك الإقليمي.", 2) | Out-Null

# 10) "We look forward to seeing you there!" -> Arabic (both occurrences)
$d.Content.Find.Execute("We look forward to seeing you there!", $true, $false, $false, $false, $false, $true, 1, $false, "نتطلع إلى رؤيتك هناك!", 2) | Out-Null

# 11) "If you have any questions, please contact us via " -> Arabic
$d.Content.Find.Execute("If you have any questions, please contact us via ", $true, $false, $false, $false, $false, $true, 1, $false, "إذا كانت لديك أي أسئلة، فاتصل بنا:  ", 2) | Out-Null

# 12) "live chat" -> "الدردشة الحية"
$d.Content.Find.Execute("live chat", $true, $false, $false, $false, $false, $true, 1, $false, "الدردشة الحية", 2) | Out-Null

# 13) Comment text "choose either one" -> "اختر أيًا منهما"
#     Direct Range.Text assignment is the reliable path for comment bodies.
$d.Comments.Item(1).Range.Text = "اختر أيًا منهما"

